# complete previous typo fix & knit to ppt
# "discharge" -> "discharge_dttm" (3 occurrences) and tidy up a duplicated
# "from in" -> "from" wording, on slides 27 and 29.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 27 ("Exercise"): Content Placeholder 2, 2nd bullet:
#   "Display the discharge vector from in cchic."
#   -> "Display the discharge_dttm vector from cchic."
# ---------------------------------------------------------------------
$s27 = $p.Slides.Item(27)
$tr27 = $s27.Shapes.Item(2).TextFrame.TextRange

$full27 = $tr27.Text
$idx = $full27.IndexOf("discharge")
$chars = $tr27.Characters($idx + 1, 9)
$chars.Text = "discharge_dttm"

$full27b = $tr27.Text
$idx2 = $full27b.IndexOf(" vector from in ")
$chars2 = $tr27.Characters($idx2 + 1, 16)
$chars2.Text = " vector from "

# ---------------------------------------------------------------------
# Slide 29 ("Answer 2"): Content Placeholder 2
#   1st paragraph: "Display the discharge vector in cchic."
#   -> "Display the discharge_dttm vector in cchic."
#   2nd paragraph: "cchic$discharge" -> "cchic$discharge_dttm"
# ---------------------------------------------------------------------
$s29 = $p.Slides.Item(29)
$tr29 = $s29.Shapes.Item(2).TextFrame.TextRange

$full29 = $tr29.Text
$idx3 = $full29.IndexOf("discharge")
$chars3 = $tr29.Characters($idx3 + 1, 9)
$chars3.Text = "discharge_dttm"

$full29b = $tr29.Text
$idx4 = $full29b.IndexOf("discharge", $idx3 + 14)
$chars4 = $tr29.Characters($idx4 + 1, 9)
$chars4.Text = "discharge_dttm"
